# Expand abbreviated "JV" level values to "Junior Varsity" across the sheet.
# The level_1 (col L) and level_2 (col M) columns use "JV" as shorthand for
# "Junior Varsity"; replace every exact-match occurrence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$used.Replace("JV", "Junior Varsity", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
